$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Formula = '''69.653.40'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Formula = '''  +0.52%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Formula = '''3.707.58'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Formula = '''  +0.82%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Formula = '''  +0.06%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Formula = '''673.35'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Formula = '''  -1.39%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Formula = '''162.04'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Formula = '''  +2.58%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Formula = '''  -0.02%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Formula = '''0.499'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Formula = '''  +1.05%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('E9').Formula = '''  +0.85%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Formula = '''7.11'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Formula = '''  +2.03%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('E11').Formula = '''  +1.92%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('E12').Formula = '''  +1.48%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Formula = '''32.89'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Formula = '''  +2.18%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Formula = '''3.714.96'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Formula = '''  +1.01%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Formula = '''69.679.56'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Formula = '''  +0.55%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('E16').Formula = '''  +1.82%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Formula = '''16.31'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Formula = '''  +2.47%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Formula = '''6.51'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Formula = '''  +1.97%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Formula = '''474.62'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Formula = '''  +0.80%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('E20').Formula = '''  -1.61%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Formula = '''0.655'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Formula = '''  +0.95%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Formula = '''80.42'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Formula = '''  +0.55%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Formula = '''3.855.04'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E24').Formula = '''  +5.50%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('E25').Formula = '''  -0.02%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Formula = '''11.02'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Formula = '''  +1.05%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Formula = '''9.17'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Formula = '''  +0.69%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Formula = '''2.70'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Formula = '''  -0.19%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('E29').Formula = '''  -0.74%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('B30').Formula = '''Kaspa'
$ws.Range('B30').Style = 'Normal'
$ws.Range('C30').Formula = '''https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('C30').Style = 'Normal'
$ws.Range('D30').Formula = '''0.169'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Formula = '''  +7.41%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('B31').Formula = '''ImmutableX'
$ws.Range('B31').Style = 'Normal'
$ws.Range('C31').Formula = '''https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('C31').Style = 'Normal'
$ws.Range('D31').Formula = '''2.02'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Formula = '''  +1.53%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Formula = '''6.60'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Formula = '''  +0.62%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Formula = '''26.97'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Formula = '''  +0.53%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Formula = '''0.998'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Formula = '''  -0.26%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Formula = '''3.697.54'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Formula = '''  +1.23%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('E36').Formula = '''  +4.47%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Formula = '''6.13'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Formula = '''  +1.50%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('D39').Formula = '''2.27'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Formula = '''  +2.14%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Formula = '''0.999'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Formula = '''  -0.09%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Formula = '''0.0916'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Formula = '''  +1.63%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Formula = '''173.89'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Formula = '''  +3.82%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('E43').Formula = '''  +0.20%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Formula = '''47.09'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E45').Formula = '''  +2.55%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Formula = '''0.000281'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Formula = '''  +1.08%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('E47').Formula = '''  +2.50%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('B48').Formula = '''SuiNetwork'
$ws.Range('B48').Style = 'Normal'
$ws.Range('C48').Formula = '''https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('C48').Style = 'Normal'
$ws.Range('D48').Formula = '''1.10'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Formula = '''  -0.43%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('B49').Formula = '''InjectiveProtocol'
$ws.Range('B49').Style = 'Normal'
$ws.Range('C49').Formula = '''https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('C49').Style = 'Normal'
$ws.Range('D49').Formula = '''27.67'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Formula = '''  +2.44%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('E50').Formula = '''  +1.83%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('E51').Formula = '''  +0.89%  '
$ws.Range('E51').Style = 'Normal'
